# Split the run "01-02/07 :" into three runs "01-02" / "-03" / "/07 :"
# (same formatting throughout) by inserting "-03" right after "01-02".
#
# The COM-interop runtime silently coalesces adjacent runs that share
# identical run properties, so a plain InsertAfter() would just grow the
# existing run instead of producing the three separate <w:r> elements the
# diff calls for. Toggling a character property on the freshly inserted
# text (on, then back off) is enough to make the engine keep it as its own
# run while leaving the final formatting untouched.

$d = $word.ActiveDocument

$matchRange = $d.Content
$found = $matchRange.Find.Execute("01-02/07", $false, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '01-02/07' in the document"
}

$splitPoint = $matchRange.Start + 5   # length of "01-02"

$insertionRange = $d.Range($splitPoint, $splitPoint)
$insertionRange.InsertAfter("-03")

# Re-grab the freshly inserted text as its own Range and force a run break
# by flipping Bold on then off again.
$newRange = $d.Range($splitPoint, $splitPoint + 3)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0
